$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Midpoint Clamp (CSV PWM)")

$ws.Range("G2").Formula = "=H2/2-1"
$ws.Range("D4").Select()
